$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.282.56"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.404.82"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.60"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.96"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.403.71"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.68"
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "3.983.52"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.70"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").Value = "3.402.13"
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("D18").Value = "61.261.38"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.93"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.87"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.36"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.39"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "3.527.24"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.553"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.18"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("E28").Value = "  -4.41%  "
$ws.Range("E29").Value = "  +8.30%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.45"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.09"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.49"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +2.25%  "
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.87"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.81"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0771"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.23"
$ws.Range("E41").Value = "  +6.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.778"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.97"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "2.534.94"
$ws.Range("E48").Value = "  +8.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.78"
$ws.Range("E49").Value = "  +5.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.79"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.38"
$ws.Range("E51").Value = "  +1.57%  "
